$wb = $excel.ActiveWorkbook

# --- Sheet: stuff-descriptor (add "Scroll of Warding" item row) ---
$ws1 = $wb.Worksheets.Item("stuff-descriptor")
$ws1.Range("A25").Value = "WardScroll"
$ws1.Range("B25").Value = "scroll"
$ws1.Range("C25").Value = "Scroll of Warding"
$ws1.Range("D25").Value = "Cast a spell of warding on yourself."
$ws1.Range("E25").Value = "#12cbe3"
$ws1.Range("L25").Value = 5
$ws1.Range("N25").Value = 0
$ws1.Range("M25").Value = 0
$ws1.Range("N25").Copy()
$ws1.Range("M25").PasteSpecial(-4122)

# --- Sheet: item-chances (add drop chances for the new scroll) ---
$ws3 = $wb.Worksheets.Item("item-chances")
$ws3.Range("A13").Value = 3
$ws3.Range("B13").Value = "WardScroll"
$ws3.Range("C13").Value = 8
$ws3.Range("A14").Value = 4
$ws3.Range("B14").Value = "WardScroll"
$ws3.Range("C14").Value = 20

$f1 = "'stuff-descriptor'!`$A`$2:`$A`$1048576"
$ws3.Range("B2:B1012").Validation.Delete()
$ws3.Range("B2:B12").Validation.Add(3, 1, 3, $f1, "0")
$ws3.Range("B15:B1012").Validation.Add(3, 1, 3, $f1, "0")

# --- Selection / active sheet bookkeeping ---
$ws1.Range("D27").Select()
$ws3.Activate()
$ws3.Range("C15").Select()
